$d = $word.ActiveDocument

# --- 1. Remove the _GoBack bookmark from its original location (after "Mwawasi") ---
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# --- 2. Flappy-bird game: merge "Tasks:..." + "Test-Group:..." paragraphs into one sentence ---
$cr = [char]13
$search = "Tasks: Starting and closing the game, playing the game, saving milestones in the game, buying powerups in the game.${cr}Test-Group: 100 – 200 test players."
$replace = "The tests would include starting and closing the game, playing the game, saving milestones in the game, buying powerups in the game. The test would involve about a hundred test players to make sure the game was easy to learn."
$d.Content.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null

# --- 3. Phone manufacture system: change title text ---
$d.Content.Find.Execute("Development of a knowledge based system for medical prescription", $true, $false, $false, $false, $false, $true, 1, $false, "Development of a phone manufacture system", 2) | Out-Null

# --- 4. Phone manufacture system: merge "Tasks:..." + "Test:..." + "Test-Group:..." paragraphs ---
$search2 = "Tasks: Adding information to the system, deriving information from the system, editing and deleting information from the system.${cr}Test: checking if the correct prescription is given when certain symptoms are displayed.${cr}Test-Group: five to ten specialised doctors."
$replace2 = "The tests would include the safety of machinery use, the quality of the product, the speed of production, the efficiency to current systems. The test would involve about a hundred phone production experts."
$d.Content.Find.Execute($search2, $true, $false, $false, $false, $false, $true, 1, $false, $replace2, 2) | Out-Null

# --- 5. Personal-information risk #2: rewrite sentence ---
$d.Content.Find.Execute("Users may not understand the conditions and consequences of providing their personal information. Make sure to explain the terms and conditions to the user.", $true, $false, $false, $false, $false, $true, 1, $false, "The users may not be familiar with where and how their information will be used. Inform the users of how and where their information will be used.", 2) | Out-Null

# --- 6. Personal-information risk #3: rewrite sentence ---
$d.Content.Find.Execute("Laws governing the provision of personal details may be different in different places. Be sure to adhere to the laws governing the provision of personal details.", $true, $false, $false, $false, $false, $true, 1, $false, "Institutions have different ethical guidelines for the use of personal information for research purposes. Make sure to get approval from the ethics department before doing your research.", 2) | Out-Null

# --- 7. Storm et al. paragraph: remove the artificial page-break split around "walking" ---
$d.Content.Find.Execute("indoor walking, descending 24 steps, outdoor walking", $true, $false, $false, $false, $false, $true, 1, $false, "indoor walking, descending 24 steps, outdoor walking", 2) | Out-Null

# --- 8. Monitor test conditions: append more detail, keeping lightGray highlight ---
$d.Content.Find.Execute("The monitors were all worn in the same places", $true, $false, $false, $false, $false, $true, 1, $false, "The monitors were all worn in the same places by every individual to ensure that the readings were the same.", 2) | Out-Null

$d.Content.Find.Execute("Users were not allowed to make any sharp turns during the tests", $true, $false, $false, $false, $false, $true, 1, $false, "Users were not allowed to make any sharp turns during the tests to avoid any mistakes of readings.", 2) | Out-Null

$d.Content.Find.Execute("The monitors started the measurements at the same time.", $true, $false, $false, $false, $false, $true, 1, $false, "The monitors started the measurements at the same time to make sure that the monitor readings were allowed the same duration.", 2) | Out-Null

$d.Content.Find.Execute("The independent variables were the Gender, age, Weight, Height and BMI", $true, $false, $false, $false, $false, $true, 1, $false, "The independent variables were the Gender, age, Weight, Height and BMI because these variables cannot change.", 2) | Out-Null

# --- 9. Dependent variables: append more detail, then delete the trailing two paragraphs ---
$d.Content.Find.Execute("The dependent variables were: The duration, accuracy of the step detectors", $true, $false, $false, $false, $false, $true, 1, $false, "The dependent variables were: The duration, accuracy of the step detectors are dependent on the step count as well as the independent variables.", 2) | Out-Null

# Remove the (now orphaned) empty paragraph + "This is because..." paragraph that followed.
$search3 = "are dependent on the step count as well as the independent variables.${cr}${cr} This is because the walking speed determines the period the experiment took. On the other hand, the accuracy of the step detector would be dependent on the walking conditions that the experiment was done in."
$replace3 = "are dependent on the step count as well as the independent variables."
$d.Content.Find.Execute($search3, $true, $false, $false, $false, $false, $true, 1, $false, $replace3, 2) | Out-Null

# --- 10. Re-add the _GoBack bookmark at the end of the dependent-variables paragraph ---
$target = $d.Content.Find
$rng = $d.Content
$rng.Find.Execute("are dependent on the step count as well as the independent variables.") | Out-Null
$endRng = $rng.Duplicate
$endRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endRng) | Out-Null
